$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 6409.25
$ws.Range("I2").Value = 8485
$ws.Range("K2").Value = 8485
$ws.Range("M2").Value = -8372
$ws.Range("H6").Value = 9999.666999999999
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 9999.666999999999
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 29999.001
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -30223.001
$ws.Range("H11").Value = 61619.445
$ws.Range("I11").Value = 61619.445
$ws.Range("K11").Value = 61619.445
$ws.Range("M11").Value = -61479.445
$ws.Range("H15").Value = 1044.8772
$ws.Range("I15").Value = 1044.8772
$ws.Range("K15").Value = 3134.6316
$ws.Range("M15").Value = -2965.6316
$ws.Range("H18").Value = 2962.6667
$ws.Range("I18").Value = 2962.6667
$ws.Range("K18").Value = 2962.6667
$ws.Range("M18").Value = -2678.6667
$ws.Range("H57").Value = 77999.336
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 300000
$ws.Range("N57").Value = -300998
$ws.Range("H62").Value = 1355000.5
$ws.Range("J62").Value = 90000
$ws.Range("L62").Value = 90000
$ws.Range("N62").Value = -91248
$ws.Range("H65").Value = 1355000.5
$ws.Range("J65").Value = 90000
$ws.Range("L65").Value = 450000
$ws.Range("N65").Value = -456240
$ws.Range("H70").Value = 512435.66
$ws.Range("I70").Value = 1021783.3
$ws.Range("K70").Value = 3065349.9
$ws.Range("M70").Value = -3065079.9
$ws.Range("H73").Value = 512435.66
$ws.Range("I73").Value = 1021783.3
$ws.Range("K73").Value = 3065349.9
$ws.Range("M73").Value = -3064413.9
$ws.Range("H96").Value = 3070.2222
$ws.Range("I96").Value = 304.57144
$ws.Range("J96").Value = 12750
$ws.Range("K96").Value = 913.71432
$ws.Range("L96").Value = 38250
$ws.Range("M96").Value = 459.28568
$ws.Range("N96").Value = -40996
$ws.Range("H98").Value = 516.3182
$ws.Range("I98").Value = 516.3182
$ws.Range("K98").Value = 516.3182
$ws.Range("M98").Value = 981.6818
$ws.Range("H101").Value = 450
$ws.Range("I101").Value = 700
$ws.Range("K101").Value = 2100
$ws.Range("M101").Value = -478
$ws.Range("H104").Value = 2383.1667
$ws.Range("I104").Value = 1266.6666
$ws.Range("J104").Value = 3499.6667
$ws.Range("K104").Value = 3799.9998
$ws.Range("L104").Value = 10499.0001
$ws.Range("M104").Value = -2052.9998
$ws.Range("N104").Value = -13993.0001
$ws.Range("H107").Value = 574.55
$ws.Range("I107").Value = 586.625
$ws.Range("K107").Value = 586.625
$ws.Range("M107").Value = 1333.375
$ws.Range("H111").Value = 41774
$ws.Range("I111").Value = 23782.5
$ws.Range("K111").Value = 71347.5
$ws.Range("M111").Value = -68280.5
$ws.Range("H113").Value = 66670704
$ws.Range("J113").Value = 5284.1665
$ws.Range("L113").Value = 5284.1665
$ws.Range("N113").Value = -11792.1665
$ws.Range("H116").Value = 36527860
$ws.Range("I116").Value = 25105666
$ws.Range("J116").Value = 55564850
$ws.Range("K116").Value = 25105666
$ws.Range("L116").Value = 55564850
$ws.Range("M116").Value = -25102224
$ws.Range("N116").Value = -55571734
$ws.Range("H122").Value = 516.3182
$ws.Range("I122").Value = 516.3182
$ws.Range("K122").Value = 1548.9546
$ws.Range("M122").Value = 901.0454
$ws.Range("H132").Value = 3436.5342
$ws.Range("I132").Value = 3247.611
$ws.Range("J132").Value = 3973.4736
$ws.Range("K132").Value = 9742.832999999999
$ws.Range("L132").Value = 11920.4208
$ws.Range("M132").Value = -7212.832999999999
$ws.Range("N132").Value = -16980.4208
$ws.Range("H135").Value = 62500696
$ws.Range("I135").Value = 62500696
$ws.Range("K135").Value = 562506264
$ws.Range("M135").Value = -562503729
$ws.Range("H137").Value = 3874.49
$ws.Range("I137").Value = 3133
$ws.Range("J137").Value = 3897.4226
$ws.Range("K137").Value = 9399
$ws.Range("L137").Value = 11692.2678
$ws.Range("M137").Value = -6849
$ws.Range("N137").Value = -16792.2678
$ws.Range("H138").Value = 5222.28
$ws.Range("I138").Value = 4079.3
$ws.Range("J138").Value = 5508.025
$ws.Range("K138").Value = 12237.9
$ws.Range("L138").Value = 16524.075
$ws.Range("M138").Value = -7097.900000000001
$ws.Range("N138").Value = -26804.075
$ws.Range("H141").Value = 1289.5
$ws.Range("I141").Value = 1343.8889
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 4031.6667
$ws.Range("L141").Value = 2400
$ws.Range("M141").Value = 1148.3333
$ws.Range("N141").Value = -12760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 174434.9
$ws.Range("I32").Value = 187190.94
$ws.Range("J32").Value = 2228.5
$ws.Range("K32").Value = 187190.94
$ws.Range("L32").Value = 2228.5
$ws.Range("M32").Value = -186903.94
$ws.Range("N32").Value = -2802.5
$ws.Range("H45").Value = 3457.6
$ws.Range("I45").Value = 3239.5715
$ws.Range("J45").Value = 3966.3333
$ws.Range("K45").Value = 3239.5715
$ws.Range("L45").Value = 3966.3333
$ws.Range("M45").Value = -2862.5715
$ws.Range("N45").Value = -4720.3333
$ws.Range("H61").Value = 2901.5908
$ws.Range("I61").Value = 2136.0588
$ws.Range("J61").Value = 5504.4
$ws.Range("K61").Value = 2136.0588
$ws.Range("L61").Value = 5504.4
$ws.Range("M61").Value = -1924.0588
$ws.Range("N61").Value = -5928.4
$ws.Range("H63").Value = 80011336
$ws.Range("J63").Value = 20016400
$ws.Range("L63").Value = 20016400
$ws.Range("N63").Value = -20017772
$ws.Range("H66").Value = 80011336
$ws.Range("J66").Value = 20016400
$ws.Range("L66").Value = 100082000
$ws.Range("N66").Value = -100088864
$ws.Range("H74").Value = 6974.304
$ws.Range("I74").Value = 6340.1577
$ws.Range("K74").Value = 6340.1577
$ws.Range("M74").Value = -5466.1577
$ws.Range("H77").Value = 6974.304
$ws.Range("I77").Value = 6340.1577
$ws.Range("K77").Value = 31700.7885
$ws.Range("M77").Value = -27332.7885
$ws.Range("H102").Value = 29858.143
$ws.Range("I102").Value = 21802
$ws.Range("J102").Value = 49998.5
$ws.Range("K102").Value = 21802
$ws.Range("L102").Value = 49998.5
$ws.Range("M102").Value = -20180
$ws.Range("N102").Value = -53242.5
$ws.Range("H110").Value = 83347290
$ws.Range("I110").Value = 100001544
$ws.Range("K110").Value = 100001544
$ws.Range("M110").Value = -99999499
$ws.Range("H122").Value = 12822965
$ws.Range("I122").Value = 16668529
$ws.Range("K122").Value = 50005587
$ws.Range("M122").Value = -50003137
$ws.Range("H131").Value = 53665.668
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H132").Value = 25003166
$ws.Range("I132").Value = 28574514
$ws.Range("K132").Value = 85723542
$ws.Range("M132").Value = -85721012
$ws.Range("H136").Value = 2901.5908
$ws.Range("I136").Value = 2136.0588
$ws.Range("J136").Value = 5504.4
$ws.Range("K136").Value = 6408.176399999999
$ws.Range("L136").Value = 16513.2
$ws.Range("M136").Value = -3858.176399999999
$ws.Range("N136").Value = -21613.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13591.272
$ws.Range("J20").Value = 5849.75
$ws.Range("L20").Value = 5849.75
$ws.Range("N20").Value = -6343.75
$ws.Range("H25").Value = 2651.25
$ws.Range("I25").Value = 439.8
$ws.Range("J25").Value = 6337
$ws.Range("K25").Value = 439.8
$ws.Range("L25").Value = 6337
$ws.Range("M25").Value = -204.8
$ws.Range("N25").Value = -6807
$ws.Range("H43").Value = 200000
$ws.Range("J43").Value = 200000
$ws.Range("L43").Value = 200000
$ws.Range("N43").Value = -200362
$ws.Range("H86").Value = 31252334
$ws.Range("I86").Value = 55557388
$ws.Range("J86").Value = 2978.1428
$ws.Range("K86").Value = 55557388
$ws.Range("L86").Value = 2978.1428
$ws.Range("M86").Value = -55556265
$ws.Range("N86").Value = -5224.1428
$ws.Range("H89").Value = 31252334
$ws.Range("I89").Value = 55557388
$ws.Range("J89").Value = 2978.1428
$ws.Range("K89").Value = 277786940
$ws.Range("L89").Value = 14890.714
$ws.Range("M89").Value = -277781324
$ws.Range("N89").Value = -26122.714
$ws.Range("H94").Value = 27783054
$ws.Range("I94").Value = 41670830
$ws.Range("J94").Value = 7500
$ws.Range("K94").Value = 41670830
$ws.Range("L94").Value = 7500
$ws.Range("M94").Value = -41670379
$ws.Range("N94").Value = -8402
$ws.Range("H99").Value = 1940
$ws.Range("I99").Value = 1955
$ws.Range("J99").Value = 1910
$ws.Range("K99").Value = 1955
$ws.Range("L99").Value = 1910
$ws.Range("M99").Value = -457
$ws.Range("N99").Value = -4906
$ws.Range("H105").Value = 1999.8
$ws.Range("I105").Value = 1999.8
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1999.8
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -252.8
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 23838100
$ws.Range("I107").Value = 14178.154
$ws.Range("J107").Value = 62551972
$ws.Range("K107").Value = 14178.154
$ws.Range("L107").Value = 62551972
$ws.Range("M107").Value = -12258.154
$ws.Range("N107").Value = -62555812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1155
$ws.Range("I16").Value = 1155
$ws.Range("K16").Value = 1155
$ws.Range("M16").Value = -868
$ws.Range("H22").Value = 2362.2942
$ws.Range("I22").Value = 1456.1428
$ws.Range("K22").Value = 1456.1428
$ws.Range("M22").Value = -1106.1428
$ws.Range("H31").Value = 3854.9453
$ws.Range("I31").Value = 979.625
$ws.Range("J31").Value = 4208.8306
$ws.Range("K31").Value = 979.625
$ws.Range("L31").Value = 4208.8306
$ws.Range("M31").Value = -684.625
$ws.Range("N31").Value = -4798.8306
$ws.Range("H34").Value = 3854.9453
$ws.Range("I34").Value = 979.625
$ws.Range("J34").Value = 4208.8306
$ws.Range("K34").Value = 979.625
$ws.Range("L34").Value = 4208.8306
$ws.Range("M34").Value = -777.625
$ws.Range("N34").Value = -4612.8306
$ws.Range("H55").Value = 18000
$ws.Range("I55").Value = 18000
$ws.Range("K55").Value = 18000
$ws.Range("M55").Value = -17685
$ws.Range("H58").Value = 359764.47
$ws.Range("I58").Value = 1326.8572
$ws.Range("K58").Value = 1326.8572
$ws.Range("M58").Value = -1123.8572
$ws.Range("H99").Value = 2448.3572
$ws.Range("I99").Value = 2156.2222
$ws.Range("J99").Value = 2974.2
$ws.Range("K99").Value = 2156.2222
$ws.Range("L99").Value = 2974.2
$ws.Range("M99").Value = -658.2222000000002
$ws.Range("N99").Value = -5970.2
$ws.Range("H105").Value = 1846
$ws.Range("I105").Value = 1855.4
$ws.Range("J105").Value = 1799
$ws.Range("K105").Value = 1855.4
$ws.Range("L105").Value = 1799
$ws.Range("M105").Value = -108.4000000000001
$ws.Range("N105").Value = -5293
$ws.Range("H113").Value = 1155
$ws.Range("I113").Value = 1155
$ws.Range("K113").Value = 1155
$ws.Range("M113").Value = 1015
$ws.Range("H122").Value = 3998.75
$ws.Range("I122").Value = 3998.75
$ws.Range("K122").Value = 11996.25
$ws.Range("M122").Value = -9546.25
$ws.Range("H126").Value = 2448.3572
$ws.Range("I126").Value = 2156.2222
$ws.Range("J126").Value = 2974.2
$ws.Range("K126").Value = 6468.6666
$ws.Range("L126").Value = 8922.599999999999
$ws.Range("M126").Value = -3998.6666
$ws.Range("N126").Value = -13862.6
$ws.Range("H132").Value = 1178406
$ws.Range("I132").Value = 770607.9
$ws.Range("K132").Value = 2311823.7
$ws.Range("M132").Value = -2309293.7
$ws.Range("H134").Value = 3197.3635
$ws.Range("I134").Value = 2617.6924
$ws.Range("K134").Value = 7853.0772
$ws.Range("M134").Value = -5318.0772
$ws.Range("H136").Value = 359764.47
$ws.Range("I136").Value = 1326.8572
$ws.Range("K136").Value = 3980.5716
$ws.Range("M136").Value = -1430.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 495.1111
$ws.Range("I13").Value = 87
$ws.Range("K13").Value = 261
$ws.Range("M13").Value = -93
$ws.Range("H56").Value = 5500
$ws.Range("I56").Value = 5500
$ws.Range("K56").Value = 5500
$ws.Range("M56").Value = -4970
$ws.Range("H57").Value = 2952.5
$ws.Range("I57").Value = 2952.5
$ws.Range("K57").Value = 8857.5
$ws.Range("M57").Value = -8298.5
$ws.Range("H68").Value = 1268.7693
$ws.Range("I68").Value = 1082.8334
$ws.Range("K68").Value = 3248.5002
$ws.Range("M68").Value = -2437.5002
$ws.Range("H71").Value = 1268.7693
$ws.Range("I71").Value = 1082.8334
$ws.Range("K71").Value = 9745.500599999999
$ws.Range("M71").Value = -5689.500599999999
$ws.Range("H86").Value = 867
$ws.Range("I86").Value = 745.6667
$ws.Range("J86").Value = 1049
$ws.Range("K86").Value = 2237.0001
$ws.Range("L86").Value = 3147
$ws.Range("M86").Value = -1051.0001
$ws.Range("N86").Value = -5519
$ws.Range("H89").Value = 867
$ws.Range("I89").Value = 745.6667
$ws.Range("J89").Value = 1049
$ws.Range("K89").Value = 6711.0003
$ws.Range("L89").Value = 9441
$ws.Range("M89").Value = -783.0002999999997
$ws.Range("N89").Value = -21297
$ws.Range("H109").Value = 1119.8
$ws.Range("I109").Value = 999
$ws.Range("J109").Value = 1150
$ws.Range("K109").Value = 2997
$ws.Range("L109").Value = 3450
$ws.Range("M109").Value = -1957
$ws.Range("N109").Value = -5530
$ws.Range("H120").Value = 36224.066
$ws.Range("I120").Value = 25999.334
$ws.Range("J120").Value = 38780.25
$ws.Range("K120").Value = 77998.00199999999
$ws.Range("L120").Value = 116340.75
$ws.Range("M120").Value = -73160.00199999999
$ws.Range("N120").Value = -126016.75
$ws.Range("H121").Value = 138705.38
$ws.Range("J121").Value = 221398.8
$ws.Range("L121").Value = 664196.3999999999
$ws.Range("N121").Value = -666816.3999999999
$ws.Range("H136").Value = 4132.375
$ws.Range("I136").Value = 4437
$ws.Range("K136").Value = 13311
$ws.Range("M136").Value = -8211
$ws.Range("H137").Value = 2488.1667
$ws.Range("J137").Value = 2988.889
$ws.Range("L137").Value = 8966.667000000001
$ws.Range("N137").Value = -19166.667
$ws.Range("H138").Value = 4074491.8
$ws.Range("I138").Value = 5455580
$ws.Range("J138").Value = 276499.5
$ws.Range("K138").Value = 16366740
$ws.Range("L138").Value = 829498.5
$ws.Range("M138").Value = -16361600
$ws.Range("N138").Value = -839778.5
$ws.Range("H139").Value = 1451187.8
$ws.Range("I139").Value = 2084776.2
$ws.Range("K139").Value = 6254328.6
$ws.Range("M139").Value = -6249188.6
$ws.Range("H141").Value = 1398.8334
$ws.Range("I141").Value = 1398.8334
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4196.5002
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 983.4997999999996
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -326
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -600
$ws.Range("H80").Value = 3378.3125
$ws.Range("I80").Value = 3412.375
$ws.Range("K80").Value = 3412.375
$ws.Range("M80").Value = -2414.375
$ws.Range("H83").Value = 3378.3125
$ws.Range("I83").Value = 3412.375
$ws.Range("K83").Value = 17061.875
$ws.Range("M83").Value = -12069.875
$ws.Range("H97").Value = 1889.2916
$ws.Range("I97").Value = 1953.5
$ws.Range("J97").Value = 1696.6666
$ws.Range("K97").Value = 1953.5
$ws.Range("L97").Value = 1696.6666
$ws.Range("M97").Value = -1457.5
$ws.Range("N97").Value = -2688.6666
$ws.Range("H102").Value = 2118.5454
$ws.Range("I102").Value = 717.5
$ws.Range("J102").Value = 3799.8
$ws.Range("K102").Value = 717.5
$ws.Range("L102").Value = 3799.8
$ws.Range("M102").Value = 904.5
$ws.Range("N102").Value = -7043.8
$ws.Range("H107").Value = 2793
$ws.Range("J107").Value = 3932.6667
$ws.Range("L107").Value = 3932.6667
$ws.Range("N107").Value = -7772.6667
$ws.Range("H113").Value = 3168.9546
$ws.Range("I113").Value = 2008
$ws.Range("J113").Value = 3972.6924
$ws.Range("K113").Value = 2008
$ws.Range("L113").Value = 3972.6924
$ws.Range("M113").Value = 162
$ws.Range("N113").Value = -8312.6924
$ws.Range("H122").Value = 100002480
$ws.Range("I122").Value = 2473.3333
$ws.Range("K122").Value = 7419.999899999999
$ws.Range("M122").Value = -4969.999899999999
$ws.Range("H126").Value = 9677.117
$ws.Range("I126").Value = 15501.5
$ws.Range("J126").Value = 4499.8887
$ws.Range("K126").Value = 46504.5
$ws.Range("L126").Value = 13499.6661
$ws.Range("M126").Value = -44034.5
$ws.Range("N126").Value = -18439.6661
$ws.Range("H132").Value = 178514.48
$ws.Range("I132").Value = 246773.31
$ws.Range("J132").Value = 3601.25
$ws.Range("K132").Value = 740319.9299999999
$ws.Range("L132").Value = 10803.75
$ws.Range("M132").Value = -737789.9299999999
$ws.Range("N132").Value = -15863.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3588040.8
$ws.Range("I22").Value = 4160
$ws.Range("J22").Value = 8067891.5
$ws.Range("K22").Value = 4160
$ws.Range("L22").Value = 8067891.5
$ws.Range("M22").Value = -3865
$ws.Range("N22").Value = -8068481.5
$ws.Range("H27").Value = 3588040.8
$ws.Range("I27").Value = 4160
$ws.Range("J27").Value = 8067891.5
$ws.Range("K27").Value = 4160
$ws.Range("L27").Value = 8067891.5
$ws.Range("M27").Value = -4053
$ws.Range("N27").Value = -8068105.5
$ws.Range("H40").Value = 2828.353
$ws.Range("I40").Value = 2930.5386
$ws.Range("J40").Value = 2496.25
$ws.Range("K40").Value = 2930.5386
$ws.Range("L40").Value = 2496.25
$ws.Range("M40").Value = -2794.5386
$ws.Range("N40").Value = -2768.25
$ws.Range("H55").Value = 403.4091
$ws.Range("I55").Value = 417.33334
$ws.Range("J55").Value = 340.75
$ws.Range("K55").Value = 417.33334
$ws.Range("L55").Value = 340.75
$ws.Range("M55").Value = -244.33334
$ws.Range("N55").Value = -686.75
$ws.Range("H58").Value = 250001550
$ws.Range("I58").Value = 2074
$ws.Range("K58").Value = 2074
$ws.Range("M58").Value = -1814
$ws.Range("H61").Value = 269898.44
$ws.Range("I61").Value = 319316.9
$ws.Range("J61").Value = 6333.3335
$ws.Range("K61").Value = 319316.9
$ws.Range("L61").Value = 6333.3335
$ws.Range("M61").Value = -319114.9
$ws.Range("N61").Value = -6737.3335
$ws.Range("H93").Value = 6633.3335
$ws.Range("I93").Value = 7450
$ws.Range("K93").Value = 7450
$ws.Range("M93").Value = -6202
$ws.Range("H113").Value = 269898.44
$ws.Range("I113").Value = 319316.9
$ws.Range("J113").Value = 6333.3335
$ws.Range("K113").Value = 319316.9
$ws.Range("L113").Value = 6333.3335
$ws.Range("M113").Value = -317146.9
$ws.Range("N113").Value = -10673.3335
$ws.Range("H122").Value = 3721.4285
$ws.Range("I122").Value = 2810
$ws.Range("K122").Value = 8430
$ws.Range("M122").Value = -5980
$ws.Range("H132").Value = 6458.7393
$ws.Range("I132").Value = 3561.6667
$ws.Range("K132").Value = 10685.0001
$ws.Range("M132").Value = -8155.000100000001
$ws.Range("H136").Value = 5039.8164
$ws.Range("I136").Value = 4528.85
$ws.Range("K136").Value = 13586.55
$ws.Range("M136").Value = -11036.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 196993.33
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H80").Value = 19500
$ws.Range("J80").Value = 19500
$ws.Range("L80").Value = 19500
$ws.Range("N80").Value = -21496
$ws.Range("H81").Value = 18190036
$ws.Range("I81").Value = 4950
$ws.Range("K81").Value = 9900
$ws.Range("M81").Value = -8839
$ws.Range("H82").Value = 7500
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 19500
$ws.Range("J83").Value = 19500
$ws.Range("L83").Value = 58500
$ws.Range("N83").Value = -68484
$ws.Range("H84").Value = 18190036
$ws.Range("I84").Value = 4950
$ws.Range("K84").Value = 49500
$ws.Range("M84").Value = -44196
$ws.Range("H85").Value = 7500
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 715.5
$ws.Range("I107").Value = 498.25
$ws.Range("K107").Value = 1494.75
$ws.Range("M107").Value = 425.25
$ws.Range("H113").Value = 919
$ws.Range("I113").Value = 967.4545000000001
$ws.Range("J113").Value = 741.3333
$ws.Range("K113").Value = 2902.3635
$ws.Range("L113").Value = 2223.9999
$ws.Range("M113").Value = -732.3635000000004
$ws.Range("N113").Value = -6563.9999
$ws.Range("H122").Value = 2105.2856
$ws.Range("I122").Value = 2036.7273
$ws.Range("K122").Value = 6110.1819
$ws.Range("M122").Value = -3660.1819
$ws.Range("H126").Value = 1352.8572
$ws.Range("I126").Value = 1228.3334
$ws.Range("K126").Value = 3685.0002
$ws.Range("M126").Value = -1215.0002
$ws.Range("H132").Value = 545225
$ws.Range("I132").Value = 1004577.4
$ws.Range("K132").Value = 3013732.2
$ws.Range("M132").Value = -3011202.2
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280
